# Generate Report for Handback
#
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns for the 0a0d7e78-... source file row on both the zh-cn and de-de
# handback-status sheets, reflecting a newly (re)generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 12:20:22"
$wsZhCn.Range("K2").Value = "2016-08-31 12:21:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-31 12:20:34"
$wsDeDe.Range("K2").Value = "2016-08-31 12:21:56"
